$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "time ago" relative text in column B for the rows that changed
$ws.Range("B2").Value = "31 minutes ago"
$ws.Range("B20").Value = "17 hours ago"
$ws.Range("B21").Value = "17 hours ago"
$ws.Range("B22").Value = "17 hours ago"
$ws.Range("B49").Value = "18 hours ago"

# Update the crawl_time column (F) for every data row (2-200) to the new timestamp
for ($r = 2; $r -le 200; $r++) {
    $ws.Cells.Item($r, 6).Value = "2025-08-15 12:25:01"
}
